$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the existing trailer rows (old 90/96/98) down to their new positions
# --- (107, 113, 115) and clear the old row contents first to avoid collisions
# --- with the newly inserted experiment rows 84-99.
$ws.Range("B90:P90").ClearContents()
$ws.Range("B96:P96").ClearContents()
$ws.Range("B98:P98").ClearContents()

# --- Fill in new experiment rows 83-99 (row 83 also gains extra columns C:P) ---
# Row 83
$ws.Range("B83").Value = 82
$ws.Range("C83").Value = 'Stat values from NSAA\AD w/ seq_len=10 (w/ scaling seq_overlap) to perform overall NSAA score regression'
$ws.Range("D83").Value = 'python comp_stat_vals.py NSAA AD all --split_size=1 '
$ws.Range("E83").Value = 'python ft_sel_red.py NSAA AD all pca --num_features=30 --no_normalize '
$ws.Range("F83").Value = 'python rnn.py NSAA AD all overall --seq_len=10 --seq_overlap=0.9'
$ws.Range("G83").Value = 'Mean Squared Error = 0.0093, Mean Absolute Error = 0.066, Root Mean Squared Error = 0.0967, R^2 Score = 0.9998'
$ws.Range("H83").Value = 'X shape = (13240, 10, 30)'
$ws.Range("I83").Value = 'Y shape = (13240,)'
$ws.Range("J83").Value = 'Test ratio = 0.2'
$ws.Range("K83").Value = 'Sequence length = 10'
$ws.Range("L83").Value = 'Features length = 30'
$ws.Range("M83").Value = 'Num epochs = 300'
$ws.Range("N83").Value = 'Num LSTM units per layer = 128'
$ws.Range("O83").Value = 'Num hidden layers = 2'
$ws.Range("P83").Value = 'Learning rate = 0.001'

# Row 84
$ws.Range("B84").Value = 83
$ws.Range("C84").Value = 'Stat values from NSAA\AD w/ seq_len=7 (w/ scaling seq_overlap) to perform overall NSAA score regression'
$ws.Range("D84").Value = 'python comp_stat_vals.py NSAA AD all --split_size=1 '
$ws.Range("E84").Value = 'python ft_sel_red.py NSAA AD all pca --num_features=30 --no_normalize '
$ws.Range("F84").Value = 'python rnn.py NSAA AD all overall --seq_len=7 --seq_overlap=0.857143'
$ws.Range("G84").Value = 'Mean Squared Error = 0.0458, Mean Absolute Error = 0.1516, Root Mean Squared Error = 0.2139, R^2 Score = 0.999'
$ws.Range("H84").Value = 'X shape = (13265, 7, 30)'
$ws.Range("I84").Value = 'Y shape = (13265,)'
$ws.Range("J84").Value = 'Test ratio = 0.2'
$ws.Range("K84").Value = 'Sequence length = 7'
$ws.Range("L84").Value = 'Features length = 30'
$ws.Range("M84").Value = 'Num epochs = 300'
$ws.Range("N84").Value = 'Num LSTM units per layer = 128'
$ws.Range("O84").Value = 'Num hidden layers = 2'
$ws.Range("P84").Value = 'Learning rate = 0.001'

# Row 85
$ws.Range("B85").Value = 84
$ws.Range("C85").Value = 'Stat values from NSAA\AD w/ seq_len=5 (w/ scaling seq_overlap) to perform overall NSAA score regression'
$ws.Range("D85").Value = 'python comp_stat_vals.py NSAA AD all --split_size=1 '
$ws.Range("E85").Value = 'python ft_sel_red.py NSAA AD all pca --num_features=30 --no_normalize '
$ws.Range("F85").Value = 'python rnn.py NSAA AD all overall --seq_len=5 --seq_overlap=0.8'
$ws.Range("G85").Value = 'Mean Squared Error = 0.1063, Mean Absolute Error = 0.1993, Root Mean Squared Error = 0.326, R^2 Score = 0.9975'
$ws.Range("H85").Value = 'X shape = (13305, 5, 30)'
$ws.Range("I85").Value = 'Y shape = (13305,)'
$ws.Range("J85").Value = 'Test ratio = 0.2'
$ws.Range("K85").Value = 'Sequence length = 5'
$ws.Range("L85").Value = 'Features length = 30'
$ws.Range("M85").Value = 'Num epochs = 300'
$ws.Range("N85").Value = 'Num LSTM units per layer = 128'
$ws.Range("O85").Value = 'Num hidden layers = 2'
$ws.Range("P85").Value = 'Learning rate = 0.001'

# Row 86
$ws.Range("B86").Value = 85
$ws.Range("C86").Value = 'Stat values from NSAA\AD w/ seq_len=3 (w/ scaling seq_overlap) to perform overall NSAA score regression'
$ws.Range("D86").Value = 'python comp_stat_vals.py NSAA AD all --split_size=1 '
$ws.Range("E86").Value = 'python ft_sel_red.py NSAA AD all pca --num_features=30 --no_normalize '
$ws.Range("F86").Value = 'python rnn.py NSAA AD all overall --seq_len=3 --seq_overlap=0.67'
$ws.Range("G86").Value = 'Mean Squared Error = 0.124, Mean Absolute Error = 0.1937, Root Mean Squared Error = 0.3522, R^2 Score = 0.9972'
$ws.Range("H86").Value = 'X shape = (13453, 3, 30)'
$ws.Range("I86").Value = 'Y shape = (13453,)'
$ws.Range("J86").Value = 'Test ratio = 0.2'
$ws.Range("K86").Value = 'Sequence length = 3'
$ws.Range("L86").Value = 'Features length = 30'
$ws.Range("M86").Value = 'Num epochs = 300'
$ws.Range("N86").Value = 'Num LSTM units per layer = 128'
$ws.Range("O86").Value = 'Num hidden layers = 2'
$ws.Range("P86").Value = 'Learning rate = 0.001'

# Row 87
$ws.Range("B87").Value = 86
$ws.Range("C87").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression (diff from #14 as 27% more data)'
$ws.Range("D87").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E87").Value = '(Not used)'
$ws.Range("F87").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=60'
$ws.Range("G87").Value = 'Mean Squared Error = 4.4579, Mean Absolute Error = 1.2406, Root Mean Squared Error = 2.1114, R^2 Score = 0.899'
$ws.Range("H87").Value = 'X shape = (13365, 60, 66)'
$ws.Range("I87").Value = 'Y shape = (13365,)'
$ws.Range("J87").Value = 'Test ratio = 0.2'
$ws.Range("K87").Value = 'Sequence length = 60'
$ws.Range("L87").Value = 'Features length = 66'
$ws.Range("M87").Value = 'Num epochs = 20'
$ws.Range("N87").Value = 'Num LSTM units per layer = 128'
$ws.Range("O87").Value = 'Num hidden layers = 2'
$ws.Range("P87").Value = 'Learning rate = 0.001'

# Row 88
$ws.Range("B88").Value = 87
$ws.Range("C88").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 50% more sequence length w/ corresponding overlap (no discard_prop)'
$ws.Range("D88").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E88").Value = '(Not used)'
$ws.Range("F88").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=90 --seq_overlap=0.333'
$ws.Range("G88").Value = 'Mean Squared Error = 4.9932, Mean Absolute Error = 1.2672, Root Mean Squared Error = 2.2346, R^2 Score = 0.8864'
$ws.Range("H88").Value = 'X shape = (13314, 90, 66)'
$ws.Range("I88").Value = 'Y shape = (13314,)'
$ws.Range("J88").Value = 'Test ratio = 0.2'
$ws.Range("K88").Value = 'Sequence length = 90'
$ws.Range("L88").Value = 'Features length = 66'
$ws.Range("M88").Value = 'Num epochs = 20'
$ws.Range("N88").Value = 'Num LSTM units per layer = 128'
$ws.Range("O88").Value = 'Num hidden layers = 2'
$ws.Range("P88").Value = 'Learning rate = 0.001'

# Row 89
$ws.Range("B89").Value = 88
$ws.Range("C89").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 100% more sequence length w/ corresponding overlap (no discard_prop)'
$ws.Range("D89").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E89").Value = '(Not used)'
$ws.Range("F89").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=120 --seq_overlap=0.5'
$ws.Range("G89").Value = 'Mean Squared Error = 4.2065, Mean Absolute Error = 1.1576, Root Mean Squared Error = 2.051, R^2 Score = 0.9103'
$ws.Range("H89").Value = 'X shape = (13315, 120, 66)'
$ws.Range("I89").Value = 'Y shape = (13315,)'
$ws.Range("J89").Value = 'Test ratio = 0.2'
$ws.Range("K89").Value = 'Sequence length = 120'
$ws.Range("L89").Value = 'Features length = 66'
$ws.Range("M89").Value = 'Num epochs = 20'
$ws.Range("N89").Value = 'Num LSTM units per layer = 128'
$ws.Range("O89").Value = 'Num hidden layers = 2'
$ws.Range("P89").Value = 'Learning rate = 0.001'

# Row 90
$ws.Range("B90").Value = 89
$ws.Range("C90").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 200% more sequence length w/ corresponding overlap (no discard_prop)'
$ws.Range("D90").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E90").Value = '(Not used)'
$ws.Range("F90").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=180 --seq_overlap=0.667'
$ws.Range("G90").Value = 'Mean Squared Error = 2.6475, Mean Absolute Error = 0.9896, Root Mean Squared Error = 1.6271, R^2 Score = 0.9427'
$ws.Range("H90").Value = 'X shape = (13313, 180, 66)'
$ws.Range("I90").Value = 'Y shape = (13313,)'
$ws.Range("J90").Value = 'Test ratio = 0.2'
$ws.Range("K90").Value = 'Sequence length = 180'
$ws.Range("L90").Value = 'Features length = 66'
$ws.Range("M90").Value = 'Num epochs = 20'
$ws.Range("N90").Value = 'Num LSTM units per layer = 128'
$ws.Range("O90").Value = 'Num hidden layers = 2'
$ws.Range("P90").Value = 'Learning rate = 0.001'

# Row 91
$ws.Range("B91").Value = 90
$ws.Range("C91").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 50% more sequence length w/ corresponding overlap w/ corresponding discard_prop'
$ws.Range("D91").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E91").Value = '(Not used)'
$ws.Range("F91").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=90 --seq_overlap=0.333 --discard_prop=0.333'
$ws.Range("G91").Value = 'Mean Squared Error = 3.3633, Mean Absolute Error = 1.2131, Root Mean Squared Error = 1.8339, R^2 Score = 0.9257'
$ws.Range("H91").Value = 'X shape = (13314, 60, 66)'
$ws.Range("I91").Value = 'Y shape = (13314,)'
$ws.Range("J91").Value = 'Test ratio = 0.2'
$ws.Range("K91").Value = 'Sequence length = 60'
$ws.Range("L91").Value = 'Features length = 66'
$ws.Range("M91").Value = 'Num epochs = 20'
$ws.Range("N91").Value = 'Num LSTM units per layer = 128'
$ws.Range("O91").Value = 'Num hidden layers = 2'
$ws.Range("P91").Value = 'Learning rate = 0.001'

# Row 92
$ws.Range("B92").Value = 91
$ws.Range("C92").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 100% more sequence length w/ corresponding overlap w/ corresponding discard_prop'
$ws.Range("D92").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E92").Value = '(Not used)'
$ws.Range("F92").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=120 --seq_overlap=0.5 --discard_prop=0.5'
$ws.Range("G92").Value = 'Mean Squared Error = 3.4289, Mean Absolute Error = 1.1326, Root Mean Squared Error = 1.8517, R^2 Score = 0.9261'
$ws.Range("H92").Value = 'X shape = (13315, 60, 66)'
$ws.Range("I92").Value = 'Y shape = (13315,)'
$ws.Range("J92").Value = 'Test ratio = 0.2'
$ws.Range("K92").Value = 'Sequence length = 60'
$ws.Range("L92").Value = 'Features length = 66'
$ws.Range("M92").Value = 'Num epochs = 20'
$ws.Range("N92").Value = 'Num LSTM units per layer = 128'
$ws.Range("O92").Value = 'Num hidden layers = 2'
$ws.Range("P92").Value = 'Learning rate = 0.001'

# Row 93
$ws.Range("B93").Value = 92
$ws.Range("C93").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 200% more sequence length w/ corresponding overlap w/ corresponding discard_prop'
$ws.Range("D93").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E93").Value = '(Not used)'
$ws.Range("F93").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=180 --seq_overlap=0.667 --discard_prop=0.667'
$ws.Range("G93").Value = 'Mean Squared Error = 3.1308, Mean Absolute Error = 1.0699, Root Mean Squared Error = 1.7694, R^2 Score = 0.9326'
$ws.Range("H93").Value = 'X shape = (13313, 60, 66)'
$ws.Range("I93").Value = 'Y shape = (13313,)'
$ws.Range("J93").Value = 'Test ratio = 0.2'
$ws.Range("K93").Value = 'Sequence length = 60'
$ws.Range("L93").Value = 'Features length = 66'
$ws.Range("M93").Value = 'Num epochs = 20'
$ws.Range("N93").Value = 'Num LSTM units per layer = 128'
$ws.Range("O93").Value = 'Num hidden layers = 2'
$ws.Range("P93").Value = 'Learning rate = 0.001'

# Row 94
$ws.Range("B94").Value = 93
$ws.Range("C94").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 400% more sequence length w/ corresponding overlap w/ corresponding discard_prop'
$ws.Range("D94").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E94").Value = '(Not used)'
$ws.Range("F94").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=300 --seq_overlap=0.8 --discard_prop=0.8'
$ws.Range("G94").Value = 'Mean Squared Error = 2.4002, Mean Absolute Error = 0.8787, Root Mean Squared Error = 1.5493, R^2 Score = 0.9466'
$ws.Range("H94").Value = 'X shape = (13240, 60, 66)'
$ws.Range("I94").Value = 'Y shape = (13240,)'
$ws.Range("J94").Value = 'Test ratio = 0.2'
$ws.Range("K94").Value = 'Sequence length = 60'
$ws.Range("L94").Value = 'Features length = 66'
$ws.Range("M94").Value = 'Num epochs = 20'
$ws.Range("N94").Value = 'Num LSTM units per layer = 128'
$ws.Range("O94").Value = 'Num hidden layers = 2'
$ws.Range("P94").Value = 'Learning rate = 0.001'

# Row 95
$ws.Range("B95").Value = 94
$ws.Range("C95").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 900% more sequence length w/ corresponding overlap w/ corresponding discard_prop'
$ws.Range("D95").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E95").Value = '(Not used)'
$ws.Range("F95").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'
$ws.Range("G95").Value = 'Mean Squared Error = 1.0565, Mean Absolute Error = 0.5997, Root Mean Squared Error = 1.0279, R^2 Score = 0.9767'
$ws.Range("H95").Value = 'X shape = (13051, 60, 66)'
$ws.Range("I95").Value = 'Y shape = (13051,)'
$ws.Range("J95").Value = 'Test ratio = 0.2'
$ws.Range("K95").Value = 'Sequence length = 60'
$ws.Range("L95").Value = 'Features length = 66'
$ws.Range("M95").Value = 'Num epochs = 20'
$ws.Range("N95").Value = 'Num LSTM units per layer = 128'
$ws.Range("O95").Value = 'Num hidden layers = 2'
$ws.Range("P95").Value = 'Learning rate = 0.001'

# Row 96
$ws.Range("B96").Value = 95
$ws.Range("C96").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 1900% more sequence length w/ corresponding overlap w/ corresponding discard_prop'
$ws.Range("D96").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E96").Value = '(Not used)'
$ws.Range("F96").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=1200 --seq_overlap=0.95 --discard_prop=0.95'
$ws.Range("G96").Value = 'Mean Squared Error = 1.4799, Mean Absolute Error = 0.657, Root Mean Squared Error = 1.2165, R^2 Score = 0.967'
$ws.Range("H96").Value = 'X shape = (12436, 60, 66)'
$ws.Range("I96").Value = 'Y shape = (12436,)'
$ws.Range("J96").Value = 'Test ratio = 0.2'
$ws.Range("K96").Value = 'Sequence length = 60'
$ws.Range("L96").Value = 'Features length = 66'
$ws.Range("M96").Value = 'Num epochs = 20'
$ws.Range("N96").Value = 'Num LSTM units per layer = 128'
$ws.Range("O96").Value = 'Num hidden layers = 2'
$ws.Range("P96").Value = 'Learning rate = 0.001'

# Row 97
$ws.Range("B97").Value = 96
$ws.Range("C97").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 3900% more sequence length w/ corresponding overlap w/ corresponding discard_prop'
$ws.Range("D97").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E97").Value = '(Not used)'
$ws.Range("F97").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=2400 --seq_overlap=0.975 --discard_prop=0.975'
$ws.Range("G97").Value = 'Mean Squared Error = 0.3534, Mean Absolute Error = 0.406, Root Mean Squared Error = 0.5945, R^2 Score = 0.9917'
$ws.Range("H97").Value = 'X shape = (11530, 60, 66)'
$ws.Range("I97").Value = 'Y shape = (11530,)'
$ws.Range("J97").Value = 'Test ratio = 0.2'
$ws.Range("K97").Value = 'Sequence length = 60'
$ws.Range("L97").Value = 'Features length = 66'
$ws.Range("M97").Value = 'Num epochs = 20'
$ws.Range("N97").Value = 'Num LSTM units per layer = 128'
$ws.Range("O97").Value = 'Num hidden layers = 2'
$ws.Range("P97").Value = 'Learning rate = 0.001'

# Row 98
$ws.Range("B98").Value = 97
$ws.Range("C98").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 7900% more sequence length w/ corresponding overlap w/ corresponding discard_prop'
$ws.Range("D98").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E98").Value = '(Not used)'
$ws.Range("F98").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=4800 --seq_overlap=0.9875 --discard_prop=0.9875'
$ws.Range("G98").Value = 'Mean Squared Error = 0.3486, Mean Absolute Error = 0.4562, Root Mean Squared Error = 0.5905, R^2 Score = 0.9926'
$ws.Range("H98").Value = 'X shape = (10130, 60, 66)'
$ws.Range("I98").Value = 'Y shape = (10130,)'
$ws.Range("J98").Value = 'Test ratio = 0.2'
$ws.Range("K98").Value = 'Sequence length = 60'
$ws.Range("L98").Value = 'Features length = 66'
$ws.Range("M98").Value = 'Num epochs = 20'
$ws.Range("N98").Value = 'Num LSTM units per layer = 128'
$ws.Range("O98").Value = 'Num hidden layers = 2'
$ws.Range("P98").Value = 'Learning rate = 0.001'

# Row 99
$ws.Range("B99").Value = 98
$ws.Range("C99").Value = 'Raw jointAngle values from NSAA\AD to perform overall NSAA score regression w/ 15900% more sequence length w/ corresponding overlap w/ corresponding discard_prop'
$ws.Range("D99").Value = 'python ext_raw_measures.py NSAA all all'
$ws.Range("E99").Value = '(Not used)'
$ws.Range("F99").Value = 'python rnn.py NSAA jointAngle all overall --seq_len=9600 --seq_overlap=0.99375 --discard_prop=0.99375'
$ws.Range("G99").Value = 'Mean Squared Error = 0.7334, Mean Absolute Error = 0.6984, Root Mean Squared Error = 0.8564, R^2 Score = 0.9803'
$ws.Range("H99").Value = 'X shape = (7481, 60, 66)'
$ws.Range("I99").Value = 'Y shape = (7481,)'
$ws.Range("J99").Value = 'Test ratio = 0.2'
$ws.Range("K99").Value = 'Sequence length = 60'
$ws.Range("L99").Value = 'Features length = 66'
$ws.Range("M99").Value = 'Num epochs = 20'
$ws.Range("N99").Value = 'Num LSTM units per layer = 128'
$ws.Range("O99").Value = 'Num hidden layers = 2'
$ws.Range("P99").Value = 'Learning rate = 0.001'

# --- Re-create the trailer rows at their new row numbers ---
# Row 107
$ws.Range("C107").Value = 'Raw joint angles from allmatfiles to perform overall NSAA score regression'
$ws.Range("D107").Value = 'python ext_raw_measures.py allmatfiles all jointAngle'

# Row 113
$ws.Range("C113").Value = 'NOTE: received many more files (15/06), hence more available training data'

# Row 115
$ws.Range("C115").Value = '(NOT SURE IF POSSIBLE)'
$ws.Range("D115").Value = 'Single-act stat values from NSAA\AD to perform D/HC classification'
$ws.Range("E115").Value = 'python mat_act_div.py V1 all; python comp_stat_vals.py NSAA AD all --split_size=1  --single_act'

# --- Update the view state (active cell / selection) ---
$ws.Range("R18").Select()
